# Implemented CRUD minimal Api for customers
# Insert a new customer record (004 / Ujjwal) as the new second row of the
# "Customer Info" sheet/table, pushing the previous sample row (1 / User)
# down to row 3. The table and sheet dimension grow from A1:E2 to A1:E3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing data row (row 2) down to row 3, leaving row 2 blank
# for the new customer. This preserves the old row's values/types exactly.
$ws.Rows(2).Insert()

# Code ("004") and Phone ("8787656789") look numeric, so force them to be
# stored as text (matching the source data) via a Text number format.
$ws.Range("A2").NumberFormat = "@"
$ws.Cells.Item(2, 1).Value = "004"
$ws.Cells.Item(2, 2).Value = "Ujjwal"
$ws.Cells.Item(2, 3).Value = "ujwjal@in.com"
$ws.Range("D2").NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "8787656789"
$ws.Cells.Item(2, 5).Value = 10

# Grow the table (and with it the autoFilter / sheet dimension) to cover
# the newly added row.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:E3"))
